$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "runs/balls/fours/sixes" columns store numeric-looking data as TEXT
# (Excel's "number stored as text" case). Writing a bare numeric literal via
# COM would re-type the cell as a Number, so each value is entered with a
# leading apostrophe to force/preserve text storage, matching the source data.

# Row 3
$ws.Range("C3").Value = "'25"
$ws.Range("D3").Value = "'13"
$ws.Range("E3").Value = "'2"

# Row 4
$ws.Range("C4").Value = "'48"
$ws.Range("D4").Value = "'34"
$ws.Range("E4").Value = "'4"

# Row 5
$ws.Range("C5").Value = "'43"
$ws.Range("D5").Value = "'35"
$ws.Range("E5").Value = "'4"
$ws.Range("F5").Value = "'0"

# Row 7
$ws.Range("C7").Value = "'87"
$ws.Range("D7").Value = "'53"
$ws.Range("E7").Value = "'11"
$ws.Range("F7").Value = "'1"

# Row 8
$ws.Range("C8").Value = "'22"
$ws.Range("D8").Value = "'19"
$ws.Range("E8").Value = "'4"

# Row 11
$ws.Range("C11").Value = "'58"
$ws.Range("D11").Value = "'44"
$ws.Range("E11").Value = "'6"

# Row 12
$ws.Range("C12").Value = "'72"
$ws.Range("D12").Value = "'37"
$ws.Range("E12").Value = "'1"
$ws.Range("F12").Value = "'7"

# Row 14
$ws.Range("C14").Value = "'10"
$ws.Range("D14").Value = "'9"
$ws.Range("E14").Value = "'1"
